$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; this shifts all existing rows 81..136 down to 82..137,
# copying formatting (the D-column date style) from the row being pushed down.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new weekly data point.
$ws.Range("A81").Value = 1
$ws.Range("B81").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C81").Value = "Arica y Parinacota"
$ws.Range("D81").Value = 44827
$ws.Range("E81").Value = 15
$ws.Range("F81").Value = 100112042
$ws.Range("G81").Value = "Locoto"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 160
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 22000
$ws.Range("M81").Value = 21000
$ws.Range("N81").Value = "$/caja 20 kilos"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 1050
$ws.Range("Q81").Value = 20
$ws.Range("R81").Value = "Hortaliza"
